# Deposit sheet (存款, sheet index 3): add bank / deposit_type / currency columns,
# shift owner/total into place, and append legislator metadata columns (category,
# date, legislator_name, legislator_id, source_file, index) to match the other sheets.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Header row
$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# Data rows (2-17)
# row 2: source index 155
$ws.Cells.Item(2, 1).Value = 155
$ws.Cells.Item(2, 2).Value = "合作金庫商業銀行古亭支庫"
$ws.Cells.Item(2, 3).Value = "支票存款"
$ws.Cells.Item(2, 4).Value = "新臺幣"
$ws.Cells.Item(2, 5).Value = "薛凌"
$ws.Cells.Item(2, 6).Value = 58260
$ws.Cells.Item(2, 7).Value = "deposit"
$ws.Cells.Item(2, 8).Value = "normal"
$ws.Cells.Item(2, 9).NumberFormat = "@"
$ws.Cells.Item(2, 9).Value = "2011-11-14"
$ws.Cells.Item(2, 10).Value = "薛凌"
$ws.Cells.Item(2, 11).Value = 1384
$ws.Cells.Item(2, 12).Value = "tmp7f9c1"
$ws.Cells.Item(2, 13).Value = 155

# row 3: source index 156
$ws.Cells.Item(3, 1).Value = 156
$ws.Cells.Item(3, 2).Value = "合作金庫商業銀行古亭支庫"
$ws.Cells.Item(3, 3).Value = "支票存款"
$ws.Cells.Item(3, 4).Value = "新臺幣"
$ws.Cells.Item(3, 5).Value = "薛凌"
$ws.Cells.Item(3, 6).Value = 1278251
$ws.Cells.Item(3, 7).Value = "deposit"
$ws.Cells.Item(3, 8).Value = "normal"
$ws.Cells.Item(3, 9).NumberFormat = "@"
$ws.Cells.Item(3, 9).Value = "2011-11-14"
$ws.Cells.Item(3, 10).Value = "薛凌"
$ws.Cells.Item(3, 11).Value = 1384
$ws.Cells.Item(3, 12).Value = "tmp7f9c1"
$ws.Cells.Item(3, 13).Value = 156

# row 4: source index 157
$ws.Cells.Item(4, 1).Value = 157
$ws.Cells.Item(4, 2).Value = "合作金庫商業銀行古亭支庫"
$ws.Cells.Item(4, 3).Value = "支票存款"
$ws.Cells.Item(4, 4).Value = "新臺幣"
$ws.Cells.Item(4, 5).Value = "陳勝宏"
$ws.Cells.Item(4, 6).Value = 2247864
$ws.Cells.Item(4, 7).Value = "deposit"
$ws.Cells.Item(4, 8).Value = "normal"
$ws.Cells.Item(4, 9).NumberFormat = "@"
$ws.Cells.Item(4, 9).Value = "2011-11-14"
$ws.Cells.Item(4, 10).Value = "薛凌"
$ws.Cells.Item(4, 11).Value = 1384
$ws.Cells.Item(4, 12).Value = "tmp7f9c1"
$ws.Cells.Item(4, 13).Value = 157

# row 5: source index 158
$ws.Cells.Item(5, 1).Value = 158
$ws.Cells.Item(5, 2).Value = "陽信商業銀行永和分行"
$ws.Cells.Item(5, 3).Value = "支票存款"
$ws.Cells.Item(5, 4).Value = "新臺幣"
$ws.Cells.Item(5, 5).Value = "薛凌"
$ws.Cells.Item(5, 6).Value = 552213
$ws.Cells.Item(5, 7).Value = "deposit"
$ws.Cells.Item(5, 8).Value = "normal"
$ws.Cells.Item(5, 9).NumberFormat = "@"
$ws.Cells.Item(5, 9).Value = "2011-11-14"
$ws.Cells.Item(5, 10).Value = "薛凌"
$ws.Cells.Item(5, 11).Value = 1384
$ws.Cells.Item(5, 12).Value = "tmp7f9c1"
$ws.Cells.Item(5, 13).Value = 158

# row 6: source index 159
$ws.Cells.Item(6, 1).Value = 159
$ws.Cells.Item(6, 2).Value = "陽信商業銀行永和分行"
$ws.Cells.Item(6, 3).Value = "支票存款"
$ws.Cells.Item(6, 4).Value = "新臺幣"
$ws.Cells.Item(6, 5).Value = "薛凌"
$ws.Cells.Item(6, 6).Value = 2359642
$ws.Cells.Item(6, 7).Value = "deposit"
$ws.Cells.Item(6, 8).Value = "normal"
$ws.Cells.Item(6, 9).NumberFormat = "@"
$ws.Cells.Item(6, 9).Value = "2011-11-14"
$ws.Cells.Item(6, 10).Value = "薛凌"
$ws.Cells.Item(6, 11).Value = 1384
$ws.Cells.Item(6, 12).Value = "tmp7f9c1"
$ws.Cells.Item(6, 13).Value = 159

# row 7: source index 160
$ws.Cells.Item(7, 1).Value = 160
$ws.Cells.Item(7, 2).Value = "陽信商業銀行新和分行"
$ws.Cells.Item(7, 3).Value = "支票存款"
$ws.Cells.Item(7, 4).Value = "新臺幣"
$ws.Cells.Item(7, 5).Value = "薛凌"
$ws.Cells.Item(7, 6).Value = 4716527
$ws.Cells.Item(7, 7).Value = "deposit"
$ws.Cells.Item(7, 8).Value = "normal"
$ws.Cells.Item(7, 9).NumberFormat = "@"
$ws.Cells.Item(7, 9).Value = "2011-11-14"
$ws.Cells.Item(7, 10).Value = "薛凌"
$ws.Cells.Item(7, 11).Value = 1384
$ws.Cells.Item(7, 12).Value = "tmp7f9c1"
$ws.Cells.Item(7, 13).Value = 160

# row 8: source index 161
$ws.Cells.Item(8, 1).Value = 161
$ws.Cells.Item(8, 2).Value = "陽信商業銀行石牌分行"
$ws.Cells.Item(8, 3).Value = "支票存款"
$ws.Cells.Item(8, 4).Value = "新臺幣"
$ws.Cells.Item(8, 5).Value = "陳勝宏"
$ws.Cells.Item(8, 6).Value = 97297
$ws.Cells.Item(8, 7).Value = "deposit"
$ws.Cells.Item(8, 8).Value = "normal"
$ws.Cells.Item(8, 9).NumberFormat = "@"
$ws.Cells.Item(8, 9).Value = "2011-11-14"
$ws.Cells.Item(8, 10).Value = "薛凌"
$ws.Cells.Item(8, 11).Value = 1384
$ws.Cells.Item(8, 12).Value = "tmp7f9c1"
$ws.Cells.Item(8, 13).Value = 161

# row 9: source index 162
$ws.Cells.Item(9, 1).Value = 162
$ws.Cells.Item(9, 2).Value = "陽信商業銀行石脾分行"
$ws.Cells.Item(9, 3).Value = "支票存款"
$ws.Cells.Item(9, 4).Value = "新臺幣"
$ws.Cells.Item(9, 5).Value = "陳勝宏"
$ws.Cells.Item(9, 6).Value = 20352
$ws.Cells.Item(9, 7).Value = "deposit"
$ws.Cells.Item(9, 8).Value = "normal"
$ws.Cells.Item(9, 9).NumberFormat = "@"
$ws.Cells.Item(9, 9).Value = "2011-11-14"
$ws.Cells.Item(9, 10).Value = "薛凌"
$ws.Cells.Item(9, 11).Value = 1384
$ws.Cells.Item(9, 12).Value = "tmp7f9c1"
$ws.Cells.Item(9, 13).Value = 162

# row 10: source index 163
$ws.Cells.Item(10, 1).Value = 163
$ws.Cells.Item(10, 2).Value = "陽信商業銀行社子分行"
$ws.Cells.Item(10, 3).Value = "支票存款"
$ws.Cells.Item(10, 4).Value = "新臺幣"
$ws.Cells.Item(10, 5).Value = "陳勝宏"
$ws.Cells.Item(10, 6).Value = 764
$ws.Cells.Item(10, 7).Value = "deposit"
$ws.Cells.Item(10, 8).Value = "normal"
$ws.Cells.Item(10, 9).NumberFormat = "@"
$ws.Cells.Item(10, 9).Value = "2011-11-14"
$ws.Cells.Item(10, 10).Value = "薛凌"
$ws.Cells.Item(10, 11).Value = 1384
$ws.Cells.Item(10, 12).Value = "tmp7f9c1"
$ws.Cells.Item(10, 13).Value = 163

# row 11: source index 164
$ws.Cells.Item(11, 1).Value = 164
$ws.Cells.Item(11, 2).Value = "陽信商業銀行社子分行"
$ws.Cells.Item(11, 3).Value = "支票存款"
$ws.Cells.Item(11, 4).Value = "新臺幣"
$ws.Cells.Item(11, 5).Value = "陳勝宏"
$ws.Cells.Item(11, 6).Value = 548910
$ws.Cells.Item(11, 7).Value = "deposit"
$ws.Cells.Item(11, 8).Value = "normal"
$ws.Cells.Item(11, 9).NumberFormat = "@"
$ws.Cells.Item(11, 9).Value = "2011-11-14"
$ws.Cells.Item(11, 10).Value = "薛凌"
$ws.Cells.Item(11, 11).Value = 1384
$ws.Cells.Item(11, 12).Value = "tmp7f9c1"
$ws.Cells.Item(11, 13).Value = 164

# row 12: source index 165
$ws.Cells.Item(12, 1).Value = 165
$ws.Cells.Item(12, 2).Value = "陽信商業銀行石牌分行"
$ws.Cells.Item(12, 3).Value = "支票存款"
$ws.Cells.Item(12, 4).Value = "新臺幣"
$ws.Cells.Item(12, 5).Value = "薛凌"
$ws.Cells.Item(12, 6).Value = 310074
$ws.Cells.Item(12, 7).Value = "deposit"
$ws.Cells.Item(12, 8).Value = "normal"
$ws.Cells.Item(12, 9).NumberFormat = "@"
$ws.Cells.Item(12, 9).Value = "2011-11-14"
$ws.Cells.Item(12, 10).Value = "薛凌"
$ws.Cells.Item(12, 11).Value = 1384
$ws.Cells.Item(12, 12).Value = "tmp7f9c1"
$ws.Cells.Item(12, 13).Value = 165

# row 13: source index 166
$ws.Cells.Item(13, 1).Value = 166
$ws.Cells.Item(13, 2).Value = "華泰商業銀行營業部"
$ws.Cells.Item(13, 3).Value = "支票存款"
$ws.Cells.Item(13, 4).Value = "新臺幣"
$ws.Cells.Item(13, 5).Value = "薛凌"
$ws.Cells.Item(13, 6).Value = 3024949
$ws.Cells.Item(13, 7).Value = "deposit"
$ws.Cells.Item(13, 8).Value = "normal"
$ws.Cells.Item(13, 9).NumberFormat = "@"
$ws.Cells.Item(13, 9).Value = "2011-11-14"
$ws.Cells.Item(13, 10).Value = "薛凌"
$ws.Cells.Item(13, 11).Value = 1384
$ws.Cells.Item(13, 12).Value = "tmp7f9c1"
$ws.Cells.Item(13, 13).Value = 166

# row 14: source index 167
$ws.Cells.Item(14, 1).Value = 167
$ws.Cells.Item(14, 2).Value = "板信商業銀行民生分行"
$ws.Cells.Item(14, 3).Value = "支票存款"
$ws.Cells.Item(14, 4).Value = "新臺幣"
$ws.Cells.Item(14, 5).Value = "薛凌"
$ws.Cells.Item(14, 6).Value = 15374
$ws.Cells.Item(14, 7).Value = "deposit"
$ws.Cells.Item(14, 8).Value = "normal"
$ws.Cells.Item(14, 9).NumberFormat = "@"
$ws.Cells.Item(14, 9).Value = "2011-11-14"
$ws.Cells.Item(14, 10).Value = "薛凌"
$ws.Cells.Item(14, 11).Value = 1384
$ws.Cells.Item(14, 12).Value = "tmp7f9c1"
$ws.Cells.Item(14, 13).Value = 167

# row 15: source index 168
$ws.Cells.Item(15, 1).Value = 168
$ws.Cells.Item(15, 2).Value = "國泰世華商業銀行復興分行"
$ws.Cells.Item(15, 3).Value = "支票存款"
$ws.Cells.Item(15, 4).Value = "新臺幣"
$ws.Cells.Item(15, 5).Value = "薛凌"
$ws.Cells.Item(15, 6).Value = 10315
$ws.Cells.Item(15, 7).Value = "deposit"
$ws.Cells.Item(15, 8).Value = "normal"
$ws.Cells.Item(15, 9).NumberFormat = "@"
$ws.Cells.Item(15, 9).Value = "2011-11-14"
$ws.Cells.Item(15, 10).Value = "薛凌"
$ws.Cells.Item(15, 11).Value = 1384
$ws.Cells.Item(15, 12).Value = "tmp7f9c1"
$ws.Cells.Item(15, 13).Value = 168

# row 16: source index 169
$ws.Cells.Item(16, 1).Value = 169
$ws.Cells.Item(16, 2).Value = "國泰世華商業銀行復興分行"
$ws.Cells.Item(16, 3).Value = "支票存款"
$ws.Cells.Item(16, 4).Value = "新臺幣"
$ws.Cells.Item(16, 5).Value = "薛凌"
$ws.Cells.Item(16, 6).Value = 4685
$ws.Cells.Item(16, 7).Value = "deposit"
$ws.Cells.Item(16, 8).Value = "normal"
$ws.Cells.Item(16, 9).NumberFormat = "@"
$ws.Cells.Item(16, 9).Value = "2011-11-14"
$ws.Cells.Item(16, 10).Value = "薛凌"
$ws.Cells.Item(16, 11).Value = 1384
$ws.Cells.Item(16, 12).Value = "tmp7f9c1"
$ws.Cells.Item(16, 13).Value = 169

# row 17: source index 170
$ws.Cells.Item(17, 1).Value = 170
$ws.Cells.Item(17, 2).Value = "陽信商業銀行永和分行"
$ws.Cells.Item(17, 3).Value = "支票存款"
$ws.Cells.Item(17, 4).Value = "美金"
$ws.Cells.Item(17, 5).Value = "薛凌"
$ws.Cells.Item(17, 6).Value = 56763
$ws.Cells.Item(17, 7).Value = "deposit"
$ws.Cells.Item(17, 8).Value = "normal"
$ws.Cells.Item(17, 9).NumberFormat = "@"
$ws.Cells.Item(17, 9).Value = "2011-11-14"
$ws.Cells.Item(17, 10).Value = "薛凌"
$ws.Cells.Item(17, 11).Value = 1384
$ws.Cells.Item(17, 12).Value = "tmp7f9c1"
$ws.Cells.Item(17, 13).Value = 170
